$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.850.08'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.86%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.496.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.77%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '192.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.98%  '
$ws.Range("E7").Value = '  +0.73%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.213'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.32%  '
$ws.Range("E10").Value = '  +2.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.44'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000307'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.60'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.058.20'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '616.71'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.921.68'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '12.67'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.10%  '
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.514.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.21%  '
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.990'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.82'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '105.36'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +10.73%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.83%  '
$ws.Range("E26").Value = '  +3.99%  '
$ws.Range("E27").Value = '  -0.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.89'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.28'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.16'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.27'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +11.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.63'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.44%  '
$ws.Range("E33").Value = '  -0.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '64.19'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.739.62'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '521.69'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.08'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0795'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.391'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.14%  '
$ws.Range("E41").Value = '  +1.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '36.65'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.84%  '
$ws.Range("E43").Value = '  -0.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0463'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.64%  '
$ws.Range("E45").Value = '  -2.67%  '
$ws.Range("E46").Value = '  +2.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.27%  '
$ws.Range("E48").Value = '  +0.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.75'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.27%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.09'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.52%  '
$ws.Range("E51").Value = '  -6.05%  '
